# The original paragraph reads:
#   "This will be part 2 ... done with flexbox. If this is the sort ... entitled:"
# It needs to become three paragraphs:
#   1) "This will be part 2 ... done with flexbox."
#   2) <empty paragraph>
#   3) " If this is the sort ... entitled:"   (note: leading space preserved)
#
# We do this with a single Find/Replace that injects two paragraph marks plus a
# temporary marker token for the (otherwise empty) middle paragraph, then strip
# the marker out with a second Find/Replace so the middle paragraph ends up with
# no run/text at all (matching the rest of the document's empty paragraphs).

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "flexbox. If this is the sort",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "flexbox.^p@@TEMP_MARK@@^p If this is the sort",
    2
)

$d.Content.Find.Execute(
    "@@TEMP_MARK@@",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2
)
